$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for Width and Height
$ws.Range("E1").Value = "Width"
$ws.Range("F1").Value = "Height"

# Update Local file paths (replace forward slash with backslash) and add Width/Height values
for ($r = 2; $r -le 8; $r++) {
    $d = $ws.Cells.Item($r, 4).Value2
    $d = $d -replace "21212/005", "21212\005"
    $ws.Cells.Item($r, 4).Value = $d

    $ws.Cells.Item($r, 5).Value = 900
    $ws.Cells.Item($r, 6).Value = 900
}
